$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Product Backlog")

# Text / label cells, entered in the same order the author typed them
# (so the shared-string table lands in the same sequence as the target file).
$ws.Range("C18").Value = "Main Page"
$ws.Range("C19").Value = "Doctor Sign Up Page"
$ws.Range("C20").Value = "Patient Sign Up Page"
$ws.Range("C21").Value = "Admin Sign Up Page"
$ws.Range("D18").Value = "All users"
$ws.Range("D21").Value = "Admin"
$ws.Range("D19").Value = "Doctor"
$ws.Range("D20").Value = "Patient"
$ws.Range("H19").Value = "Given that I am a doctor, when I get a job in Neighborhood Doctors and am ready to make an account, then I am able to sign up with my email, password and doctor role or another services such as Google etc."
$ws.Range("H18").Value = "Given that I am a user whoever is doctor, patient or admin, when I want to learn more about Neighborhood Doctors and create an account or log in my account, then I am able to log in to the system or choose whether I should create an account. "
$ws.Range("H20").Value = "Given that I am a patient, when I feel sick and want to make an account in Neighborhood Doctors, then I am able to sign up with my email, password and patient role or another services such as Google etc."
$ws.Range("I20").Value = "No need to make a checking function for the patient role."
$ws.Range("I18").Value = "Need to add the introduction about Neighborhood Doctors."
$ws.Range("I19").Value = "Will have a checking function whether this user is a doctor."
$ws.Range("H21").Value = "Given that I am an admin, when I get a job as a developer in Neighborhood Doctors and am ready to make an account, then I am able to sign up with my email, password and admin role or another services such as Google etc."
$ws.Range("I21").Value = "Will have a checking function whether this user is a admin."
$ws.Range("C22").Value = "Doctor Log In Page"
$ws.Range("C23").Value = "Patient Log In Page"
$ws.Range("C24").Value = "Admin Log In Page"
$ws.Range("D22").Value = "Doctor"
$ws.Range("D23").Value = "Patient"
$ws.Range("D24").Value = "Admin"
$ws.Range("H22").Value = "Given that I am a doctor, when I prepare for working in Neighborhood Doctors, then I am able to log in with my email and password or another services such as Google etc."
$ws.Range("H23").Value = "Given that I am a patient, when I want to log in to Neighborhood Doctors,then I am able to log in with my email and password or another services such as Google etc."
$ws.Range("H24").Value = "Given that I am an admin, when I am about to do the maintance job in Neighborhood Doctors, then I am able to log in with my email and password or another services such as Google etc."
$ws.Range("I22").Value = "Access the doctor information from the database."
$ws.Range("I23").Value = "Access the patient information from the database."
$ws.Range("I24").Value = "Access the admin information from the database."

# Effort / priority numeric cells
$ws.Range("E18").Value = 10
$ws.Range("F18").Value = 4
$ws.Range("E19").Value = 10
$ws.Range("F19").Value = 5
$ws.Range("E20").Value = 10
$ws.Range("F20").Value = 5
$ws.Range("E21").Value = 10
$ws.Range("F21").Value = 5
$ws.Range("E22").Value = 10
$ws.Range("F22").Value = 5
$ws.Range("E23").Value = 10
$ws.Range("F23").Value = 5
$ws.Range("E24").Value = 10
$ws.Range("F24").Value = 5

# Status cells (reuse the existing "To Do " shared string)
$ws.Range("G18").Value = "To Do "
$ws.Range("G19").Value = "To Do "
$ws.Range("G20").Value = "To Do "
$ws.Range("G21").Value = "To Do "
$ws.Range("G22").Value = "To Do "
$ws.Range("G23").Value = "To Do "
$ws.Range("G24").Value = "To Do "

# Leave the view the way the author left it: zoomed to 115%, with the
# Admin Sign Up Page row range selected.
$excel.ActiveWindow.Zoom = 115
[void]$ws.Range("C22:C24").Select()

$wb.Save()
